# Actualización 10 de Mayo
# Updates the statistics tables on "Estadisticos 2P" and "Estadisticos Final"
# with the latest Blancos / Reprobados / Aprobados / Por_Apro / Promedio figures.

$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P" (columns: C Totales, D Blancos, E Reprobados, F Aprobados, G Por_Apro, H Promedio) ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

$data2P = @(
    @{ Row = 2; D = 0; E = 9;  F = 28; G = 75.68000000000001; H = 8.300000000000001 },
    @{ Row = 3; D = 0; E = 3;  F = 33; G = 91.67;              H = 8.6 },
    @{ Row = 4; D = 0; E = 4;  F = 32; G = 88.89;              H = 8.6 },
    @{ Row = 5; D = 0; E = 12; F = 23; G = 65.70999999999999;  H = 7.2 },
    @{ Row = 6; D = 0; E = 12; F = 23; G = 65.70999999999999;  H = 7.3 }
)

foreach ($entry in $data2P) {
    $r = $entry.Row
    $ws2P.Cells.Item($r, 4).Value = $entry.D
    $ws2P.Cells.Item($r, 5).Value = $entry.E
    $ws2P.Cells.Item($r, 6).Value = $entry.F
    $ws2P.Cells.Item($r, 7).Value = $entry.G
    $ws2P.Cells.Item($r, 8).Value = $entry.H
}

# --- Sheet "Estadisticos Final" (columns: C Totales, D Blancos, E Reprobados, F Aprobados, G Por_Apro, H Promedio) ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$dataFinal = @(
    @{ Row = 2; E = 9;  F = 28; G = 75.68000000000001; H = 8.5 },
    @{ Row = 3; E = 3;  F = 33; G = 91.67;              H = 8.699999999999999 },
    @{ Row = 4; E = 4;  F = 32; G = 88.89;              H = 8.699999999999999 },
    @{ Row = 5; E = 12; F = 23; G = 65.70999999999999;  H = 7.3 },
    @{ Row = 6; E = 12; F = 23; G = 65.70999999999999;  H = 7.4 }
)

foreach ($entry in $dataFinal) {
    $r = $entry.Row
    $wsFinal.Cells.Item($r, 5).Value = $entry.E
    $wsFinal.Cells.Item($r, 6).Value = $entry.F
    $wsFinal.Cells.Item($r, 7).Value = $entry.G
    $wsFinal.Cells.Item($r, 8).Value = $entry.H
}
